# Add a "Season" column (G) derived from each sample's site code
# (P = Spring, S = Summer, W = Winter) and select the final new cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = "Season"

$seasons = @{
    2  = "Spring"
    3  = "Spring"
    4  = "Summer"
    5  = "Summer"
    6  = "Winter"
    7  = "Spring"
    8  = "Spring"
    9  = "Spring"
    10 = "Summer"
    11 = "Summer"
    12 = "Summer"
    13 = "Winter"
    14 = "Winter"
    15 = "Winter"
    16 = "Spring"
    17 = "Spring"
    18 = "Summer"
    19 = "Winter"
    20 = "Spring"
    21 = "Spring"
    22 = "Spring"
    23 = "Summer"
    24 = "Summer"
    25 = "Summer"
    26 = "Winter"
    27 = "Winter"
    28 = "Winter"
}

foreach ($row in 2..28) {
    $ws.Range("G$row").Value = $seasons[$row]
}

$ws.Range("G28").Select()
